$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The commit removes the review row for ronenchen27@gmail.com / danfogel100@gmail.com
# ("I love playing this game so much. ...") which was row 6, shifting every
# subsequent review row up by one, and appends a brand-new review row at the
# bottom (stevewonder3001@gmail.com / budoyoni@gmail.com).
# ---------------------------------------------------------------------------

# Stash a "clean" (never-hyperlinked) copy of the C/D column formatting before
# we touch anything, so we can restore it later - Hyperlinks.Add() auto-applies
# Excel's blue/underlined Hyperlink look, which the source file does not use.
$ws.Range("D10").Copy()
$ws.Range("ZZ1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 1. Delete row 6 entirely - Excel shifts rows 7:10 up to 6:9 (values, shared
#    strings and cell formatting all move with the row).
$ws.Rows("6:6").Delete()

# 2. Bring formatting for the brand-new last row (10) in line with the rest
#    of the table by cloning the formats from the row above it (row 9, which
#    is itself the old row 10 shifted up - so its look is identical to every
#    other data row).
$ws.Range("A9:F9").Copy()
$ws.Range("A10:F10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# 3. Write the new review into row 10.
$ws.Range("A10").Value2 = "com.singleton.strechy"
$ws.Range("B10").Value2 = "stretchy"
$ws.Range("C10").Value2 = "stevewonder3001@gmail.com"
$ws.Range("D10").Value2 = "budoyoni@gmail.com"
$ws.Range("E10").Value2 = "27/5/2019 15:59"
$ws.Range("F10").Value2 = "I’m very grateful for this game. My kids play it all the time. They love this car game so much. I have a lot of free time. Amazing cars game."

# 4. Row heights: the custom 13.8 height that used to belong to the deleted
#    row's successor (old row 10) now belongs to row 9; the brand new row 10
#    goes back to the sheet's default 12.8 height.
$ws.Rows(9).RowHeight = 13.8
$ws.Rows(10).RowHeight = 12.8

# 5. Hyperlinks don't automatically re-target when rows shift, so rebuild the
#    whole set from scratch to match the new row numbering.
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:snizzvered@gmail.com", "", "", "snizzvered@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:krigelron@gmail.com", "", "", "krigelron@gmail.com")

$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:jorjkluni03@gmail.com", "", "", "jorjkluni03@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:shmulmaor2@gmail.com", "", "", "shmulmaor2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:vikicrestina@gmail.com", "", "", "vikicrestina@gmail.com")

$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:sugarderryfireapp@gmail.com", "", "", "sugarderryfireapp@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:sugarderryfire@gmail.com", "", "", "sugarderryfire@gmail.com")

$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:sm6502345@gmail.com", "", "", "sm6502345@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:cybworking@gmail.com", "", "", "cybworking@gmail.com")

$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:nitanfriman@gmail.com", "", "", "nitanfriman@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:ronoren61@gmail.com", "", "", "ronoren61@gmail.com")

$ws.Hyperlinks.Add($ws.Range("C8"), "mailto:eligitel@gmail.com", "", "", "eligitel@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:ronenchen27@gmail.com", "", "", "ronenchen27@gmail.com")

$ws.Hyperlinks.Add($ws.Range("C9"), "mailto:zaittomer@gmail.com", "", "", "zaittomer@gmail.com ")

$ws.Hyperlinks.Add($ws.Range("C10"), "mailto:stevewonder3001@gmail.com", "", "", "stevewonder3001@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:budoyoni@gmail.com", "", "", "budoyoni@gmail.com")

# 6. Restore the original (non-underlined, non-blue) look across C2:D10 that
#    Hyperlinks.Add() just overwrote, using the clean format stashed in (1).
$ws.Range("ZZ1").Copy()
$ws.Range("C2:D10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("ZZ1").Clear()

Write-Host "Edit complete"
